$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Delegator_Role"), shifting it and
# everything after it one column to the right, to make room for the new
# "Delegator_mil" variable.
$ws.Columns("E").Insert()

# Populate the four rows of the newly inserted column E.
$ws.Range("E1").Value = "Delegator_mil"
$ws.Range("E2").Value = "Delegation survey response log"
$ws.Range("E3").Value = "Does this individual have a military background?"
$ws.Range("E4").Value = "yes or no"

# Update the selection to reflect where the user left off editing.
$ws.Range("G3").Select()
